# Refined Evaluation to be more exact
# - Insert two new columns (O:P) before the old "Extracted Objects" column,
#   shifting the old O..U columns to Q..W.
# - Rename the old M1/N1 headers.
# - Fill the two new columns (O1/P1 headers + O2:P6 data) with their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank columns at O:P (pushes old O..U -> Q..W, and grows
#    the used range / <dimension> automatically).
$ws.Range("O:P").Insert()

# 2) Rename the existing M1 / N1 headers.
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# 3) Header text for the two newly inserted columns.
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# 4) Populate the new column values for each data row.
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1

$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 3

$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 2

$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 4

$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
